$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.158.28"
$ws.Range("E2").Value = "  +2.83%  "
$ws.Range("D3").Value = "3.106.91"
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("E4").Value = "  -0.03%  "
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "523.51"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  +1.63%  "
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.45"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  +1.18%  "
$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = "  -0.02%  "
$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.440"
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = "  +1.00%  "
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.41"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  +1.74%  "
$ws.Range("E10").Value = "  +1.19%  "
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.385"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  +3.23%  "
$ws.Range("D12").Value = "3.639.07"
$ws.Range("E12").Value = "  +0.86%  "
$ws.Range("E13").Value = "  +1.35%  "
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.15"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  +5.23%  "
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("D16").Value = "59.119.45"
$ws.Range("E16").Value = "  +2.69%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.21"
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = "  +1.80%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.103.09"
$ws.Range("E18").Value = "  +0.63%  "
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.06"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  +0.23%  "
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.21"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  +0.26%  "
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "344.89"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  +1.64%  "
$ws.Range("E22").Value = "  -0.30%  "
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.510"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  +1.92%  "
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.04"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = "  +0.64%  "
$ws.Range("E25").Value = "  -0.97%  "
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("D27").Value = "0.0₃0936"
$ws.Range("E27").Value = "  -1.02%  "
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.85"
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = "  +5.70%  "
$ws.Range("E29").Value = "  +2.08%  "
$ws.Range("E30").Value = "  +1.95%  "
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.22"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  +3.19%  "
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.12"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  +1.44%  "
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "155.31"
$ws.Range("D33").Style = $origStyle
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.66"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  +2.31%  "
$ws.Range("E35").Value = "  +5.21%  "
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "26.98"
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = "  +2.98%  "
$ws.Range("E37").Value = "  +5.38%  "
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0691"
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = "  +1.69%  "
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.97"
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = "  +2.68%  "
$ws.Range("D40").Value = "3.148.54"
$ws.Range("E40").Value = "  +0.91%  "
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.77"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  -0.63%  "
$ws.Range("E43").Value = "  -1.04%  "
$ws.Range("E44").Value = "  +5.50%  "
$ws.Range("D45").Value = "2.294.35"
$ws.Range("E45").Value = "  +1.09%  "
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0258"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  +2.89%  "
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "20.92"
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = "  +3.07%  "
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.972"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  +1.07%  "
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.06"
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = "  +3.13%  "
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.763"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  +10.98%  "
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "262.08"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  +11.19%  "
